$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.322.71'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '2.677.66'
$ws.Range('E3').Value = '  +5.74%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''518.84'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('D6').Value = '''145.73'
$ws.Range('E6').Value = '  +1.90%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').Value = '''0.569'
$ws.Range('E8').Value = '  +1.40%  '
$ws.Range('D9').Value = '2.712.52'
$ws.Range('E9').Value = '  +7.00%  '
$ws.Range('D10').Value = '''6.26'
$ws.Range('E10').Value = '  +1.80%  '
$ws.Range('D11').Value = '''0.107'
$ws.Range('E11').Value = '  +5.59%  '
$ws.Range('D12').Value = '''0.339'
$ws.Range('E12').Value = '  +2.69%  '
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').Value = '3.144.56'
$ws.Range('E14').Value = '  +5.55%  '
$ws.Range('D15').Value = '59.271.55'
$ws.Range('D16').Value = '''21.19'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('E17').Value = '  +2.73%  '
$ws.Range('D18').Value = '2.701.06'
$ws.Range('E18').Value = '  +6.41%  '
$ws.Range('D19').Value = '''357.42'
$ws.Range('E19').Value = '  +6.99%  '
$ws.Range('D20').Value = '''4.56'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '''10.47'
$ws.Range('E21').Value = '  +4.29%  '
$ws.Range('D22').Value = '''6.23'
$ws.Range('E22').Value = '  +5.01%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = '''62.01'
$ws.Range('E24').Value = '  +3.49%  '
$ws.Range('D25').Value = '''0.423'
$ws.Range('E25').Value = '  +4.11%  '
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = '''0.991'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('D28').Value = '0.0₃0818'
$ws.Range('E28').Value = '  +4.62%  '
$ws.Range('D29').Value = '''7.26'
$ws.Range('E29').Value = '  +5.38%  '
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('D31').Value = '''6.40'
$ws.Range('E31').Value = '  +9.62%  '
$ws.Range('D32').Value = '''19.16'
$ws.Range('E32').Value = '  +3.78%  '
$ws.Range('D33').Value = '''1.59'
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('D34').Value = '''150.69'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('D35').Value = '''0.975'
$ws.Range('E35').Value = '  +4.70%  '
$ws.Range('D36').Value = '''4.05'
$ws.Range('E36').Value = '  +4.36%  '
$ws.Range('D37').Value = '''1.15'
$ws.Range('E37').Value = '  +3.82%  '
$ws.Range('D38').Value = '''36.81'
$ws.Range('E38').Value = '  +2.21%  '
$ws.Range('D39').Value = '''0.850'
$ws.Range('E39').Value = '  +3.13%  '
$ws.Range('D40').Value = '''3.74'
$ws.Range('E40').Value = '  +6.34%  '
$ws.Range('D41').Value = '''1.42'
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('D42').Value = '''282.70'
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('D43').Value = '''0.621'
$ws.Range('E43').Value = '  +3.52%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '''0.0990'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''19.99'
$ws.Range('E45').Value = '  +7.44%  '
$ws.Range('D46').Value = '''0.993'
$ws.Range('D47').Value = '''0.0533'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '''0.0232'
$ws.Range('E48').Value = '  +2.36%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '''4.76'
$ws.Range('E49').Value = '  +5.63%  '
$ws.Range('D50').Value = '2.014.26'
$ws.Range('E50').Value = '  +6.88%  '
$ws.Range('E51').Value = '  -0.03%  '
